# Apply the "working on running electronic_only" edit:
#  - Append new rows 167-181 to the opv_results sheet for the FF / JSC / VOC
#    parameters under the "electronic_only" datatype (mirrors the existing
#    calc_PCE block pattern found at rows 82-86).
#  - Update the view/selection to reflect where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("opv_results")

# Model labels used in column A, repeated for each parameter block
$models = @("RF", "BRT", "SVM", "NN", "LSTM")

# New Predict/Parameter labels added in column D (shared strings 30-32)
$parameters = @("FF", "JSC", "VOC")

# Summary statistics for the first (RF) row of each new parameter block:
# num_of_data (E), R_mean (F), R_std (G), RMSE_mean (H), RMSE_std (I)
$stats = @{
    "FF"  = @(0.25258670773465902, 0.037576020888399002, 0.121427776682782, 0.016228970514301201, 447)
    "JSC" = @(0.65633543982259201, 0.050932560971677303, 0.13704075275446201, 0.0074415413768480296, 447)
    "VOC" = @(0.66585182916198404, 0.048133234329250499, 0.063665221036705799, 0.0080769472002819205, 447)
}

$row = 167
foreach ($param in $parameters) {
    $vals = $stats[$param]
    foreach ($model in $models) {
        $ws.Cells.Item($row, 1).Value = $model
        $ws.Cells.Item($row, 2).Value = "N/A"
        $ws.Cells.Item($row, 3).Value = """electronic_only"""
        $ws.Cells.Item($row, 4).Value = $param

        if ($model -eq "RF") {
            $ws.Cells.Item($row, 5).Value = $vals[0]
            $ws.Cells.Item($row, 6).Value = $vals[1]
            $ws.Cells.Item($row, 7).Value = $vals[2]
            $ws.Cells.Item($row, 8).Value = $vals[3]
            $ws.Cells.Item($row, 9).Value = $vals[4]
        }

        $row = $row + 1
    }
}

# Reflect the author's scroll position / active selection from the session
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L77:M77").Select() | Out-Null
